# fixbug tinh chiet khau don thu no
# Update discount amounts on "Đơn thu nợ" sheet and the downstream
# totals on "Lương" that depend on them.

$wb = $excel.ActiveWorkbook

$wsDonThuNo = $wb.Worksheets.Item("Đơn thu nợ")
$wsDonThuNo.Range("S2").Value = 110000
$wsDonThuNo.Range("S4").Value = 170000

$wsLuong = $wb.Worksheets.Item("Lương")
$wsLuong.Range("B10").Value = 170000
$wsLuong.Range("B34").Value = 823571.4285714286
$wsLuong.Range("B37").Value = 923571.4285714286
